$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply an inflation adjustment (16%) to the income-based percentile figures.
# B=P10, C=P50, D=P99, I=Media_Top_0.1 are the independently computed base
# metrics; the remaining columns (E,F,G,H,J) are the differences/ratios
# derived from them, recalculated using the new inflated base figures.

function Set-Row($r, $B, $C, $D, $E, $F, $G, $H, $I, $J) {
    $ws.Range("B$r").Value = $B
    $ws.Range("C$r").Value = $C
    $ws.Range("D$r").Value = $D
    $ws.Range("E$r").Value = $E
    $ws.Range("F$r").Value = $F
    $ws.Range("G$r").Value = $G
    $ws.Range("H$r").Value = $H
    $ws.Range("I$r").Value = $I
    $ws.Range("J$r").Value = $J
}

Set-Row 2 688.562633994668 1728.62326526626 35186.1617444724 1040.060631271592 33457.53847920614 2.51048078987052 20.3550203514429 488332.0497994241 282.497672923664
Set-Row 3 688.562633994668 1728.62326526626 35186.1617444724 1040.060631271592 33457.53847920614 2.51048078987052 20.3550203514429 465216.8427860898 269.1256401170976
Set-Row 4 688.562633994668 1728.62326526626 35121.9072044382 1040.060631271592 33393.28393917194 2.51048078987052 20.31784941817752 462368.2441667155 267.4777399200959
